$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Anje Falkenrath', ['{1}{B}{R}', 'Legendary Creature — Vampire', 'Haste', '{T}, Discard a card: Draw a card.', 'Whenever you discard a card, if it has madness, untap Anje Falkenrath.', '1/3'])"
$ws.Range("A3").Value = "('Ghired, Conclave Exile', ['{2}{R}{G}{W}', 'Legendary Creature — Human Shaman', 'When Ghired, Conclave Exile enters the battlefield, create a 4/4 green Rhino creature token with trample.', 'Whenever Ghired attacks, populate. The token enters the battlefield tapped and attacking. (To populate, create a token that’s a copy of a creature token you control.)', '2/5'])"
$ws.Range("A4").Value = "('Kadena, Slinking Sorcerer', ['{1}{B}{G}{U}', 'Legendary Creature — Naga Wizard', 'The first face-down creature spell you cast each turn costs {3} less to cast.', 'Whenever a face-down creature enters the battlefield under your control, draw a card.', '3/3'])"
$ws.Range("A5").Value = "('Sevinne, the Chronoclasm', ['{2}{U}{R}{W}', 'Legendary Creature — Human Wizard', 'Prevent all damage that would be dealt to Sevinne, the Chronoclasm.', 'Whenever you cast your first instant or sorcery spell from your graveyard each turn, copy that spell. You may choose new targets for the copy.', '2/2'])"

$ws.Range("A6:A26").EntireRow.Delete()
